# Add a "STREAMS" column (H) to the "TOP ÁLBUNS" sheet (3rd sheet) with the
# total streams per album, formats it like the neighbouring header/body
# cells, widens the column, and switches the active tab from the first
# sheet to this one (selecting K4), matching the source commit.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

# Header cell H1: bold, centered, integer number format, shared string "STREAMS"
$h1 = $ws3.Cells.Item(1, 8)
$h1.Value = "STREAMS"
$h1.Font.Bold = $true
$h1.HorizontalAlignment = -4108
$h1.NumberFormat = "0"

# Body cells H2:H11: centered, integer number format, album total streams
$streams = @(8185056174, 4873205472, 9380141140, 15126897327, 8547844071, 5033709566, 3980461750, 3519407554, 3951841933, 20102823684)
for ($i = 0; $i -lt $streams.Length; $i++) {
    $cell = $ws3.Cells.Item(2 + $i, 8)
    $cell.Value = $streams[$i]
    $cell.NumberFormat = "0"
    $cell.HorizontalAlignment = -4108
}

# Widen the new column to fit its content
$ws3.Columns.Item(8).ColumnWidth = 14.75

# Move the active tab/selection from "TOP ARTISTAS" to "TOP ÁLBUNS"
$ws3.Activate()
$ws3.Range("K4").Select()
